$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date serial value from 45233 to 45243
# for rows 2-5, preserving existing cell formatting.
$ws.Range("C2").Value = 45243
$ws.Range("C3").Value = 45243
$ws.Range("C4").Value = 45243
$ws.Range("C5").Value = 45243
